$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: vacate the existing text cells (A2:D4) with unique placeholders so the
# shared-string pool drops "Bmp7"/"Acvr1"/"ECs"/"sCs" entirely (refcount -> 0, compacted).
# This lets us rebuild the pool from scratch in the exact order the data was originally
# authored in (column-major), matching upstream.
$ws.Range("A2").Value = "__tmp0__"
$ws.Range("B2").Value = "__tmp1__"
$ws.Range("C2").Value = "__tmp2__"
$ws.Range("D2").Value = "__tmp3__"
$ws.Range("A3").Value = "__tmp4__"
$ws.Range("B3").Value = "__tmp5__"
$ws.Range("C3").Value = "__tmp6__"
$ws.Range("D3").Value = "__tmp7__"
$ws.Range("A4").Value = "__tmp8__"
$ws.Range("B4").Value = "__tmp9__"
$ws.Range("C4").Value = "__tmp10__"
$ws.Range("D4").Value = "__tmp11__"

# --- Step 2: write the Sending/Ligand/Receptor/Target cluster columns column-major
# (all of column A top-to-bottom, then B, then C, then D) so new unique strings are
# appended to the shared-string table in that order: FAPs, sCs, Bmp7, Acvr1, ECs.
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "sCs"
$ws.Range("A6").Value = "sCs"
$ws.Range("A7").Value = "sCs"

$ws.Range("B2").Value = "Bmp7"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("B7").Value = "Bmp7"

$ws.Range("C2").Value = "Acvr1"
$ws.Range("C3").Value = "Acvr1"
$ws.Range("C4").Value = "Acvr1"
$ws.Range("C5").Value = "Acvr1"
$ws.Range("C6").Value = "Acvr1"
$ws.Range("C7").Value = "Acvr1"

$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"

# --- Step 3: numeric columns E:T for rows 2-7.
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.668521
$ws.Range("H2").Value = 5.005563
$ws.Range("I2").Value = 0.9677024783929865
$ws.Range("J2").Value = 0.9677024783929865
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.057757666666666
$ws.Range("N2").Value = 15.173273
$ws.Range("O2").Value = 0.173378811020062
$ws.Range("P2").Value = 0.173378811020062
$ws.Range("Q2").Value = 8.438974879744334
$ws.Range("R2").Value = 75.95077391769901
$ws.Range("S2").Value = 0.1677791051249432
$ws.Range("T2").Value = 0.1677791051249432

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.668521
$ws.Range("H3").Value = 5.005563
$ws.Range("I3").Value = 0.9677024783929865
$ws.Range("J3").Value = 0.9677024783929865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.247411
$ws.Range("N3").Value = 42.742233
$ws.Range("O3").Value = 0.4883980890531961
$ws.Range("P3").Value = 0.4883980890531961
$ws.Range("Q3").Value = 23.772104449131
$ws.Range("R3").Value = 213.948940042179
$ws.Range("S3").Value = 0.4726240412191764
$ws.Range("T3").Value = 0.4726240412191764

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.668521
$ws.Range("H4").Value = 5.005563
$ws.Range("I4").Value = 0.9677024783929865
$ws.Range("J4").Value = 0.9677024783929865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.866548666666667
$ws.Range("N4").Value = 29.599646
$ws.Range("O4").Value = 0.3382230999267418
$ws.Range("P4").Value = 0.3382230999267418
$ws.Range("Q4").Value = 16.46254364785533
$ws.Range("R4").Value = 148.162892830698
$ws.Range("S4").Value = 0.3272993320488668
$ws.Range("T4").Value = 0.3272993320488668

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05568766666666666
$ws.Range("H5").Value = 0.167063
$ws.Range("I5").Value = 0.03229752160701353
$ws.Range("J5").Value = 0.03229752160701353
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.057757666666666
$ws.Range("N5").Value = 15.173273
$ws.Range("O5").Value = 0.173378811020062
$ws.Range("P5").Value = 0.173378811020062
$ws.Range("Q5").Value = 0.2816547230221111
$ws.Range("R5").Value = 2.534892507199
$ws.Range("S5").Value = 0.005599705895118767
$ws.Range("T5").Value = 0.005599705895118768

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05568766666666666
$ws.Range("H6").Value = 0.167063
$ws.Range("I6").Value = 0.03229752160701353
$ws.Range("J6").Value = 0.03229752160701353
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.247411
$ws.Range("N6").Value = 42.742233
$ws.Range("O6").Value = 0.4883980890531961
$ws.Range("P6").Value = 0.4883980890531961
$ws.Range("Q6").Value = 0.793405074631
$ws.Range("R6").Value = 7.140645671679
$ws.Range("S6").Value = 0.01577404783401972
$ws.Range("T6").Value = 0.01577404783401972

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05568766666666666
$ws.Range("H7").Value = 0.167063
$ws.Range("I7").Value = 0.03229752160701353
$ws.Range("J7").Value = 0.03229752160701353
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.866548666666667
$ws.Range("N7").Value = 29.599646
$ws.Range("O7").Value = 0.3382230999267418
$ws.Range("P7").Value = 0.3382230999267418
$ws.Range("Q7").Value = 0.5494450732997778
$ws.Range("R7").Value = 4.945005659697999
$ws.Range("S7").Value = 0.01092376787787504
$ws.Range("T7").Value = 0.01092376787787504

